$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sample PLK9 -> PLK30 in row 10 (fieldNumber + catalogNumber columns)
$ws.Range("B10").Value = "PLK30"
$ws.Range("C10").Value = "PLK30"

# Correct the species identification for that row: Atelopus rugulosus -> Dicamptodon tenebrosus
$ws.Range("K10").Value = "Dicamptodon"
$ws.Range("L10").Value = "tenebrosus"

# Move the active selection to B11 (next empty row), matching the saved view state
$ws.Range("B11").Select()
